# EST-1338: charge INSURANCE added to charge list for capex
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The charge hierarchy list has one row per charge (Name / Parent / AtPath),
# sorted alphabetically, terminating with "OTHER" as the last entry.
# Insert a new row just above "OTHER" (row 13) for the new INSURANCE charge,
# which pushes "OTHER" down into row 14 (picking up the formatting already
# present on the old, blank row 14) while keeping the sheet's row count the
# same (it was already sized through row 14).
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = "INSURANCE"
$ws.Cells.Item(13, 2).Value = "FRANCE"
# Leading apostrophe keeps this a literal text value (matches the existing
# quote-prefixed "/FRA" cells elsewhere in the column).
$ws.Cells.Item(13, 3).Value = "'/FRA"

# The insert above shifted the old "OTHER" row (and the blank row after it)
# down by one, temporarily growing the sheet to row 15; collapse that back
# down by removing the now-superfluous extra row, restoring the original
# 14-row extent.
$ws.Rows.Item(15).Delete()

# Re-assert the row height on the (moved) "OTHER" row so it keeps its
# original custom height.
$ws.Rows.Item(14).RowHeight = 15
